$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.635.08'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.848.17'
$ws.Range("E3").Value = '  -0.01%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = "'312.56"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("E6").Value = '  +0.03%  '
$ws.Range("D7").Value = "'0.4291"
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").Value = "'44.94"
$ws.Range("E9").Value = '  +2.19%  '
$ws.Range("D10").Value = "'0.07316"
$ws.Range("E10").Value = '  +1.00%  '
$ws.Range("D11").Value = "'0.8757"
$ws.Range("E11").Value = '  -2.60%  '
$ws.Range("E12").Value = '  +0.17%  '
$ws.Range("D13").Value = '1.831.35'
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").Value = "'5.325"
$ws.Range("E14").Value = '  -0.90%  '
$ws.Range("D15").Value = "'6.523"
$ws.Range("E15").Value = '  -0.80%  '
$ws.Range("D16").Value = "'0.06918"
$ws.Range("E16").Value = '  +0.91%  '
$ws.Range("E17").Value = '  -0.10%  '
$ws.Range("D18").Value = "'79.90"
$ws.Range("E18").Value = '  +2.58%  '
$ws.Range("D19").Value = "'0.000009009"
$ws.Range("E19").Value = '  +0.75%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("D21").Value = "'15.34"
$ws.Range("E21").Value = '  -0.50%  '
$ws.Range("D22").Value = '27.645.12'
$ws.Range("E22").Value = '  -0.03%  '
$ws.Range("D23").Value = "'4.954"
$ws.Range("E23").Value = '  -0.36%  '
$ws.Range("E24").Value = '  -2.33%  '
$ws.Range("D25").Value = '2.110.87'
$ws.Range("E25").Value = '  +2.08%  '
$ws.Range("D26").Value = "'1.991"
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").Value = "'155.18"
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").Value = "'18.73"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("D29").Value = "'121.75"
$ws.Range("E29").Value = '  +9.14%  '
$ws.Range("D30").Value = "'5.285"
$ws.Range("E30").Value = '  -1.19%  '
$ws.Range("E31").Value = '  +1.85%  '
$ws.Range("D32").Value = "'0.08906"
$ws.Range("E32").Value = '  -0.51%  '
$ws.Range("D33").Value = "'0.7645"
$ws.Range("E33").Value = '  -2.93%  '
$ws.Range("D34").Value = "'2.981"
$ws.Range("E34").Value = '  +1.38%  '
$ws.Range("D35").Value = "'4.552"
$ws.Range("E35").Value = '  +0.57%  '
$ws.Range("D36").Value = "'1.104"
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("E37").Value = '  -0.72%  '
$ws.Range("D38").Value = "'1.089"
$ws.Range("E38").Value = '  -0.45%  '
$ws.Range("D39").Value = "'0.01936"
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("D40").Value = "'2.817"
$ws.Range("E40").Value = '  -5.21%  '
$ws.Range("D41").Value = "'0.5075"
$ws.Range("E41").Value = '  +0.11%  '
$ws.Range("D42").Value = "'0.1655"
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").Value = "'6.758"
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D44").Value = "'8.358"
$ws.Range("E44").Value = '  +0.34%  '
$ws.Range("D45").Value = "'0.06552"
$ws.Range("E45").Value = '  -1.19%  '
$ws.Range("D46").Value = "'10.40"
$ws.Range("E46").Value = '  +0.64%  '
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D47").Value = "'0.4685"
$ws.Range("E47").Value = '  -0.72%  '
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").Value = "'104.86"
$ws.Range("E48").Value = '  -1.75%  '
$ws.Range("E49").Value = '  -0.02%  '
$ws.Range("D50").Value = "'1.621"
$ws.Range("E50").Value = '  -1.41%  '
$ws.Range("D51").Value = "'64.47"
$ws.Range("E51").Value = '  -0.51%  '
